$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 71724; B = "Ryan da Rosa";            C = "Recursos Humanos";     D = "Viagem de negocios"; E = 8; F = 45094; G = 4737.1 },
    @{ Row = 3;  A = 37091; B = "Anthony Gabriel Fogaça";   C = "Operacoes";            D = "Viagem de negocios"; E = 8; F = 45088; G = 8031.58 },
    @{ Row = 4;  A = 35905; B = "Nathan Peixoto";           C = "Marketing";            D = "Consulta medica";    E = 2; F = 45085; G = 6753.59 },
    @{ Row = 5;  A = 11749; B = "Agatha Barros";            C = "Operacoes";            D = "Consulta medica";    E = 6; F = 45102; G = 7441.16 },
    @{ Row = 6;  A = 99086; B = "Francisco Sá";             C = "Engenharia";           D = "Viagem de negocios"; E = 7; F = 45079; G = 7109.83 },
    @{ Row = 7;  A = 85495; B = "João Lucas Pereira";       C = "Recursos Humanos";     D = "Doenca";             E = 1; F = 45099; G = 5621.17 },
    @{ Row = 8;  A = 81154; B = "Kevin Andrade";            C = "Financeiro";           D = "Doenca";             E = 3; F = 45104; G = 5852.44 },
    @{ Row = 9;  A = 98861; B = "Ana Liz Porto";            C = "TI";                   D = "Problemas pessoais"; E = 8; F = 45081; G = 5895.49 },
    @{ Row = 10; A = 41730; B = "Luna Caldeira";            C = "Marketing";            D = "Viagem de negocios"; E = 8; F = 45080; G = 3870.3 },
    @{ Row = 11; A = 56370; B = "Joaquim Viana";            C = "Atendimento ao Cliente"; D = "Consulta medica";  E = 8; F = 45095; G = 5924.06 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
